$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.414.03"
$ws.Range("E2").Value = "  +0.17%  "
$ws.Range("D3").Value = "2.595.90"
$ws.Range("E3").Value = "  -0.57%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("E5").Value = "  +0.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.45"
$ws.Range("E6").Value = "  +1.20%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("E8").Value = "  +0.38%  "
$ws.Range("D9").Value = "2.617.96"
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.49"
$ws.Range("E10").Value = "  -0.77%  "
$ws.Range("E11").Value = "  -0.88%  "
$ws.Range("E12").Value = "  +2.81%  "
$ws.Range("E13").Value = "  +0.12%  "
$ws.Range("D14").Value = "3.057.28"
$ws.Range("E14").Value = "  -0.35%  "
$ws.Range("D15").Value = "58.189.88"
$ws.Range("E15").Value = "  -0.12%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.39"
$ws.Range("E16").Value = "  -2.05%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.635.49"
$ws.Range("E17").Value = "  +0.84%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0000134"
$ws.Range("E18").Value = "  -0.80%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "340.06"
$ws.Range("E19").Value = "  +1.08%  "
$ws.Range("E20").Value = "  -0.93%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.24"
$ws.Range("E21").Value = "  -0.92%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.46"
$ws.Range("E22").Value = "  +3.50%  "
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.50"
$ws.Range("E24").Value = "  +1.02%  "
$ws.Range("E25").Value = "  +1.83%  "
$ws.Range("E26").Value = "  -1.68%  "
$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  +0.60%  "
$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").Value = "2.715.73"
$ws.Range("E28").Value = "  -0.67%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.04"
$ws.Range("E29").Value = "  -0.80%  "
$ws.Range("D30").Value = "0.0₃0745"
$ws.Range("E30").Value = "  -4.78%  "
$ws.Range("E31").Value = "  -0.10%  "
$ws.Range("E32").Value = "  -5.85%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.58"
$ws.Range("E33").Value = "  -0.27%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "18.78"
$ws.Range("E34").Value = "  +0.39%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "149.63"
$ws.Range("E35").Value = "  -0.31%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.00"
$ws.Range("E36").Value = "  -1.33%  "
$ws.Range("E37").Value = "  -3.79%  "
$ws.Range("E38").Value = "  -2.16%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.869"
$ws.Range("E39").Value = "  +2.81%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.46"
$ws.Range("E40").Value = "  +3.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.05"
$ws.Range("E41").Value = "  -0.28%  "
$ws.Range("E42").Value = "  -1.84%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.998"
$ws.Range("E43").Value = "  -0.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.608"
$ws.Range("E44").Value = "  +1.35%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "269.86"
$ws.Range("E45").Value = "  +1.37%  "
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0955"
$ws.Range("E46").Value = "  -1.23%  "
$ws.Range("B47").Value = "WhiteBITCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.67"
$ws.Range("E47").Value = "  +0.24%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.74"
$ws.Range("E48").Value = "  -1.51%  "
$ws.Range("E49").Value = "  -0.77%  "
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "1.967.74"
$ws.Range("E50").Value = "  -2.61%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "18.72"
$ws.Range("E51").Value = "  +3.30%  "
